$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2262196666666666
$ws.Range("N2").Value = 0.6786589999999999
$ws.Range("O2").Value = 0.03145179203784564
$ws.Range("P2").Value = 0.03145179203784564
$ws.Range("Q2").Value = 0.4029937467244443
$ws.Range("R2").Value = 3.626943720519999
$ws.Range("S2").Value = 0.03094893297250209
$ws.Range("T2").Value = 0.03094893297250209
$ws.Range("O3").Value = 0.9636438974901603
$ws.Range("P3").Value = 0.9636438974901604
$ws.Range("S3").Value = 0.9482369194383904
$ws.Range("T3").Value = 0.9482369194383905
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03527466666666667
$ws.Range("N4").Value = 0.105824
$ws.Range("O4").Value = 0.004904310471994002
$ws.Range("P4").Value = 0.004904310471994003
$ws.Range("Q4").Value = 0.06283923185777777
$ws.Range("R4").Value = 0.5655530867199999
$ws.Range("S4").Value = 0.004825899137684849
$ws.Range("T4").Value = 0.004825899137684849
$ws.Range("G5").Value = 0.02894466666666666
$ws.Range("H5").Value = 0.08683399999999999
$ws.Range("I5").Value = 0.01598824845142267
$ws.Range("J5").Value = 0.01598824845142267
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2262196666666666
$ws.Range("N5").Value = 0.6786589999999999
$ws.Range("O5").Value = 0.03145179203784564
$ws.Range("P5").Value = 0.03145179203784564
$ws.Range("Q5").Value = 0.006547852845111109
$ws.Range("R5").Value = 0.05893067560599999
$ws.Range("S5").Value = 0.0005028590653435535
$ws.Range("T5").Value = 0.0005028590653435535
$ws.Range("G6").Value = 0.02894466666666666
$ws.Range("H6").Value = 0.08683399999999999
$ws.Range("I6").Value = 0.01598824845142267
$ws.Range("J6").Value = 0.01598824845142267
$ws.Range("O6").Value = 0.9636438974901603
$ws.Range("P6").Value = 0.9636438974901604
$ws.Range("Q6").Value = 0.2006180896866666
$ws.Range("S6").Value = 0.01540697805176997
$ws.Range("T6").Value = 0.01540697805176997
$ws.Range("G7").Value = 0.02894466666666666
$ws.Range("H7").Value = 0.08683399999999999
$ws.Range("I7").Value = 0.01598824845142267
$ws.Range("J7").Value = 0.01598824845142267
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03527466666666667
$ws.Range("N7").Value = 0.105824
$ws.Range("O7").Value = 0.004904310471994002
$ws.Range("P7").Value = 0.004904310471994003
$ws.Range("Q7").Value = 0.001021013468444444
$ws.Range("R7").Value = 0.009189121215999999
$ws.Range("S7").Value = 0.00007841133430915411
$ws.Range("T7").Value = 0.00007841133430915412
